$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.177.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.034.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +19.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.026.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.03%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.681.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.037.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.101.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.40%  "
$ws.Range("E24").Value = "  +5.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "680.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "66.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.425"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.69%  "
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0499"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  +14.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("E51").Value = "  +2.86%  "
